$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-09-07 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-09-08 Sunday", 2) | Out-Null
$d.Content.Find.Execute("456×6=", $true, $false, $false, $false, $false, $true, 1, $false, "446×5=", 2) | Out-Null
$d.Content.Find.Execute("651×9=", $true, $false, $false, $false, $false, $true, 1, $false, "238×5=", 2) | Out-Null
$d.Content.Find.Execute("540×7=", $true, $false, $false, $false, $false, $true, 1, $false, "171×7=", 2) | Out-Null
$d.Content.Find.Execute("646×8=", $true, $false, $false, $false, $false, $true, 1, $false, "371×6=", 2) | Out-Null
$d.Content.Find.Execute("132×4=", $true, $false, $false, $false, $false, $true, 1, $false, "658×8=", 2) | Out-Null
$d.Content.Find.Execute("905×4=", $true, $false, $false, $false, $false, $true, 1, $false, "549×6=", 2) | Out-Null
$d.Content.Find.Execute("979×3=", $true, $false, $false, $false, $false, $true, 1, $false, "784×6=", 2) | Out-Null
$d.Content.Find.Execute("197×8=", $true, $false, $false, $false, $false, $true, 1, $false, "692×8=", 2) | Out-Null
$d.Content.Find.Execute("193×2=", $true, $false, $false, $false, $false, $true, 1, $false, "728×4=", 2) | Out-Null
$d.Content.Find.Execute("225×8=", $true, $false, $false, $false, $false, $true, 1, $false, "390×9=", 2) | Out-Null
$d.Content.Find.Execute("902×2=", $true, $false, $false, $false, $false, $true, 1, $false, "803×3=", 2) | Out-Null
$d.Content.Find.Execute("685×8=", $true, $false, $false, $false, $false, $true, 1, $false, "696×2=", 2) | Out-Null
$d.Content.Find.Execute("231×8=", $true, $false, $false, $false, $false, $true, 1, $false, "645×7=", 2) | Out-Null
$d.Content.Find.Execute("180×4=", $true, $false, $false, $false, $false, $true, 1, $false, "672×8=", 2) | Out-Null
$d.Content.Find.Execute("739×2=", $true, $false, $false, $false, $false, $true, 1, $false, "862×2=", 2) | Out-Null
$d.Content.Find.Execute("326×2=", $true, $false, $false, $false, $false, $true, 1, $false, "386×6=", 2) | Out-Null
$d.Content.Find.Execute("255×3=", $true, $false, $false, $false, $false, $true, 1, $false, "122×6=", 2) | Out-Null
$d.Content.Find.Execute("274×8=", $true, $false, $false, $false, $false, $true, 1, $false, "325×9=", 2) | Out-Null
$d.Content.Find.Execute("709×7=", $true, $false, $false, $false, $false, $true, 1, $false, "615×7=", 2) | Out-Null
$d.Content.Find.Execute("867×6=", $true, $false, $false, $false, $false, $true, 1, $false, "124×3=", 2) | Out-Null
$d.Content.Find.Execute("848×4=", $true, $false, $false, $false, $false, $true, 1, $false, "561×6=", 2) | Out-Null
$d.Content.Find.Execute("114×7=", $true, $false, $false, $false, $false, $true, 1, $false, "966×6=", 2) | Out-Null
$d.Content.Find.Execute("425×5=", $true, $false, $false, $false, $false, $true, 1, $false, "678×2=", 2) | Out-Null
$d.Content.Find.Execute("586×4=", $true, $false, $false, $false, $false, $true, 1, $false, "394×6=", 2) | Out-Null
$d.Content.Find.Execute("864×2=", $true, $false, $false, $false, $false, $true, 1, $false, "824×7=", 2) | Out-Null
